$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 824.75
$ws.Range("I4").Value = 824.75
$ws.Range("K4").Value = 824.75
$ws.Range("M4").Value = -710.75
$ws.Range("H51").Value = 3250
$ws.Range("I51").Value = 3250
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 3250
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 1095.6
$ws.Range("I58").Value = 333.5
$ws.Range("J58").Value = 2619.8
$ws.Range("K58").Value = 1000.5
$ws.Range("L58").Value = 7859.400000000001
$ws.Range("M58").Value = -850.5
$ws.Range("N58").Value = -8159.400000000001
$ws.Range("H61").Value = 1014.5
$ws.Range("I61").Value = 1014.5
$ws.Range("K61").Value = 3043.5
$ws.Range("M61").Value = -2871.5
$ws.Range("H76").Value = 2607466.2
$ws.Range("I76").Value = 5858405
$ws.Range("J76").Value = 6715
$ws.Range("K76").Value = 5858405
$ws.Range("L76").Value = 6715
$ws.Range("M76").Value = -5858090
$ws.Range("N76").Value = -7345
$ws.Range("H79").Value = 2607466.2
$ws.Range("I79").Value = 5858405
$ws.Range("J79").Value = 6715
$ws.Range("K79").Value = 5858405
$ws.Range("L79").Value = 6715
$ws.Range("M79").Value = -5857313
$ws.Range("N79").Value = -8899
$ws.Range("H125").Value = 1956.25
$ws.Range("I125").Value = 1900
$ws.Range("J125").Value = 1964.2858
$ws.Range("K125").Value = 17100
$ws.Range("L125").Value = 17678.5722
$ws.Range("M125").Value = -14640
$ws.Range("N125").Value = -22598.5722
$ws.Range("H132").Value = 965.2432
$ws.Range("I132").Value = 1007.9355
$ws.Range("J132").Value = 744.6667
$ws.Range("K132").Value = 3023.8065
$ws.Range("L132").Value = 2234.0001
$ws.Range("M132").Value = -493.8065000000001
$ws.Range("N132").Value = -7294.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3835.709
$ws.Range("I32").Value = 2053.8604
$ws.Range("K32").Value = 2053.8604
$ws.Range("M32").Value = -1766.8604
$ws.Range("H132").Value = 1847.7693
$ws.Range("I132").Value = 1402.2
$ws.Range("K132").Value = 4206.6
$ws.Range("M132").Value = -1676.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1449.8334
$ws.Range("I99").Value = 1049.75
$ws.Range("J99").Value = 2250
$ws.Range("K99").Value = 1049.75
$ws.Range("L99").Value = 2250
$ws.Range("M99").Value = 448.25
$ws.Range("N99").Value = -5246
$ws.Range("H119").Value = 39250
$ws.Range("J119").Value = 39250
$ws.Range("L119").Value = 39250
$ws.Range("N119").Value = -48926
$ws.Range("H134").Value = 15617.333
$ws.Range("I134").Value = 16716.867
$ws.Range("K134").Value = 50150.601
$ws.Range("M134").Value = -47615.601

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 125.85714
$ws.Range("I7").Value = 148.2
$ws.Range("J7").Value = 70
$ws.Range("K7").Value = 148.2
$ws.Range("L7").Value = 70
$ws.Range("M7").Value = -35.19999999999999
$ws.Range("N7").Value = -296
$ws.Range("H102").Value = 25000
$ws.Range("J102").Value = 25000
$ws.Range("L102").Value = 25000
$ws.Range("N102").Value = -29868
$ws.Range("H105").Value = 1208.5
$ws.Range("I105").Value = 942.8333
$ws.Range("K105").Value = 942.8333
$ws.Range("M105").Value = 804.1667

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H62").Value = 4000
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 12000
$ws.Range("N62").Value = -13372
$ws.Range("H65").Value = 4000
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 36000
$ws.Range("N65").Value = -42864
$ws.Range("H88").Value = 5641.6665
$ws.Range("J88").Value = 5971.875
$ws.Range("L88").Value = 17915.625
$ws.Range("N88").Value = -18771.625
$ws.Range("H91").Value = 5641.6665
$ws.Range("J91").Value = 5971.875
$ws.Range("L91").Value = 17915.625
$ws.Range("N91").Value = -20879.625
$ws.Range("H92").Value = 317.75
$ws.Range("I92").Value = 200
$ws.Range("J92").Value = 334.57144
$ws.Range("K92").Value = 600
$ws.Range("L92").Value = 1003.71432
$ws.Range("M92").Value = 648
$ws.Range("N92").Value = -3499.71432
$ws.Range("H103").Value = 2583.7144
$ws.Range("J103").Value = 4741.5
$ws.Range("L103").Value = 14224.5
$ws.Range("N103").Value = -15982.5
$ws.Range("H104").Value = 4409.091
$ws.Range("J104").Value = 4922.222
$ws.Range("L104").Value = 14766.666
$ws.Range("N104").Value = -20008.666
$ws.Range("H114").Value = 28575712
$ws.Range("I114").Value = 714
$ws.Range("J114").Value = 47625710
$ws.Range("K114").Value = 2142
$ws.Range("L114").Value = 142877130
$ws.Range("M114").Value = 1112
$ws.Range("N114").Value = -142883638
$ws.Range("H116").Value = 125001830
$ws.Range("J116").Value = 166668670
$ws.Range("L116").Value = 500006010
$ws.Range("N116").Value = -500012894
$ws.Range("H127").Value = 1394.3334
$ws.Range("J127").Value = 1394.3334
$ws.Range("L127").Value = 4183.0002
$ws.Range("N127").Value = -14103.0002
$ws.Range("H131").Value = 786.33
$ws.Range("I131").Value = 484.4
$ws.Range("J131").Value = 802.2210700000001
$ws.Range("K131").Value = 1453.2
$ws.Range("L131").Value = 2406.66321
$ws.Range("M131").Value = 3586.8
$ws.Range("N131").Value = -12486.66321

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("N95").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2000
$ws.Range("J81").Value = 2000
$ws.Range("L81").Value = 4000
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 2000
$ws.Range("J84").Value = 2000
$ws.Range("L84").Value = 20000
$ws.Range("N84").Value = -30608
$ws.Range("H95").Value = 42200
$ws.Range("J95").Value = 42200
$ws.Range("L95").Value = 42200
$ws.Range("N95").Value = -47692
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H136").Value = 50509636
$ws.Range("I136").Value = 92597660
$ws.Range("K136").Value = 277792980
$ws.Range("M136").Value = -277790430
